$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2 = 1; 3 = 1; 4 = 0; 5 = 0; 6 = 2; 7 = 2; 8 = 0; 9 = 2; 10 = 1;
    11 = 0;
    12 = 1; 13 = 1; 14 = 0; 15 = 0; 16 = 1; 17 = 1; 18 = 0; 19 = 0; 20 = 1;
    21 = 2; 22 = 0; 23 = 0; 24 = 1; 25 = 2; 26 = 2; 27 = 1; 28 = 1; 29 = 1;
    30 = 1; 31 = 1; 32 = 3; 33 = 2; 34 = 0; 35 = 2; 36 = 2; 37 = 2; 38 = 1;
    39 = 2; 40 = 1; 41 = 1; 42 = 1; 43 = 0; 44 = 1;
    45 = 0;
    46 = 0; 47 = 2; 48 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
